$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the text month names in column B (rows 2-13) with real date
# values (first day of each month of 2021), formatted with a custom
# number format of YYYY-MM-DD HH:MM:SS.
$years  = @(2021,2021,2021,2021,2021,2021,2021,2021,2021,2021,2021,2021)
$months = @(1,2,3,4,5,6,7,8,9,10,11,12)

for ($i = 0; $i -lt $months.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Range("B$row")
    $cell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $cell.Value = Get-Date -Year $years[$i] -Month $months[$i] -Day 1 -Hour 0 -Minute 0 -Second 0
}
